# Remove the TORNADO 400W Hand Blender product row (row 6) from the
# products table, since there are no images for this color variant in
# products.json. The row being removed also carried the document's
# "_GoBack" bookmark (left over from the last edit position inside the
# "Material: stainless Steel" text); after deleting the row, Word's
# last-edit bookmark is re-created at the (now) final empty paragraph
# that follows the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 10 (1-based) is the standalone TORNADO 400W Hand Blender row - the
# last row of the table, holding product #6's details.
$t.Rows.Item(10).Delete()

# Re-create the "_GoBack" bookmark at the empty paragraph that now
# immediately follows the table (mirrors Word's own behaviour of
# relocating _GoBack to the last edited spot).
$end = $d.Content.End
$target = $d.Range($end - 1, $end - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
